$wb = $excel.ActiveWorkbook

# Map of row -> new F-column ("想去人数") value, shared by both sheets that
# list these events ("展览" and "全部类型").
$updates = @{
    2 = 305
    3 = 41
    5 = 4518
    6 = 342
    7 = 623
    8 = 283
    9 = 701
}

# Sheet "展览": rows 2-10, F10 changes 177 -> 179
$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates[$row]
}
$ws1.Cells.Item(10, 6).Value = 179

# Sheet "全部类型": rows 2-11, F11 changes 177 -> 179 (F10 there is the
# unrelated "演出" row and stays at 18)
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates[$row]
}
$ws4.Cells.Item(11, 6).Value = 179
